$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C5) from 2023-11-13 (45243) to 2023-11-14 (45244)
$ws.Range("C2:C5").Value = 45244
